$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'63.017.30"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +3.04%  "
$c = $ws.Range("D3")
$c.Value = "'3.024.83"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +1.73%  "
$c = $ws.Range("D4")
$c.Value = "'1.00"
$c.Style = "Normal"
$ws.Range("E4").Value = "  +0.16%  "
$c = $ws.Range("D5")
$c.Value = "'595.96"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.36%  "
$c = $ws.Range("D6")
$c.Value = "'152.97"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +7.08%  "
$ws.Range("E7").Value = "  +0.04%  "
$c = $ws.Range("D8")
$c.Value = "'3.021.34"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +1.76%  "
$c = $ws.Range("D9")
$c.Value = "'0.513"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +0.18%  "
$c = $ws.Range("D10")
$c.Value = "'7.01"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +16.70%  "
$ws.Range("E11").Value = "  +1.49%  "
$c = $ws.Range("D12")
$c.Value = "'0.463"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +2.75%  "
$ws.Range("E13").Value = "  +3.13%  "
$c = $ws.Range("D14")
$c.Value = "'35.72"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +5.17%  "
$ws.Range("E15").Value = "  -0.09%  "
$c = $ws.Range("D16")
$c.Value = "'3.530.49"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +1.99%  "
$c = $ws.Range("D18")
$c.Value = "'63.023.29"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +3.11%  "
$c = $ws.Range("D19")
$c.Value = "'3.026.94"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +1.18%  "
$c = $ws.Range("D20")
$c.Value = "'449.11"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +0.75%  "
$c = $ws.Range("D21")
$c.Value = "'14.24"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +2.14%  "
$c = $ws.Range("D22")
$c.Value = "'0.696"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +2.50%  "
$ws.Range("E23").Value = "  +3.26%  "
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$c = $ws.Range("D24")
$c.Value = "'82.92"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +1.96%  "
$ws.Range("B25").Value = "RenderToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$c = $ws.Range("D25")
$c.Value = "'11.42"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +7.99%  "
$c = $ws.Range("D26")
$c.Value = "'2.30"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +5.91%  "
$c = $ws.Range("D27")
$c.Value = "'12.38"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +3.68%  "
$ws.Range("E28").Value = "  +0.03%  "
$c = $ws.Range("D29")
$c.Value = "'7.50"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +4.98%  "
$ws.Range("E31").Value = "  +0.91%  "
$c = $ws.Range("D32")
$c.Value = "'1.00"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -0.10%  "
$c = $ws.Range("D33")
$c.Value = "'27.67"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +2.23%  "
$c = $ws.Range("D34")
$c.Value = "'0.110"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +1.04%  "
$c = $ws.Range("D35")
$c.Value = "'0.0₃0875"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +8.31%  "
$ws.Range("E36").Value = "  +3.07%  "
$c = $ws.Range("D37")
$c.Value = "'5.88"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +2.43%  "
$c = $ws.Range("D38")
$c.Value = "'3.13"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +10.91%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$c = $ws.Range("D39")
$c.Value = "'2.10"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +4.36%  "
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$c = $ws.Range("D40")
$c.Value = "'0.130"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +7.61%  "
$c = $ws.Range("D41")
$c.Value = "'50.59"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +1.11%  "
$c = $ws.Range("D42")
$c.Value = "'9.03"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +1.05%  "
$c = $ws.Range("D43")
$c.Value = "'44.90"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +16.89%  "
$c = $ws.Range("D44")
$c.Value = "'0.306"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +14.81%  "
$c = $ws.Range("D45")
$c.Value = "'391.81"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +2.22%  "
$ws.Range("E46").Value = "  +3.62%  "
$c = $ws.Range("D47")
$c.Value = "'2.708.23"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +1.21%  "
$c = $ws.Range("D48")
$c.Value = "'133.93"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +2.52%  "
$c = $ws.Range("D49")
$c.Value = "'26.90"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +15.94%  "
$ws.Range("E51").Value = "  +6.81%  "
